# Add a new shared-string value ("run") to row 2 of the sheet, in cell A2,
# and move the active selection to that cell (from the previous E18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "run"
$ws.Range("A2").Select()
